$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.935.59"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "1.574.06"
$ws.Range("E3").Value = "  -2.18%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").Value = "'1.005"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "'298.48"
$ws.Range("E6").Value = "  -1.44%  "

$ws.Range("D7").Value = "'0.3733"
$ws.Range("E7").Value = "  -1.24%  "

$ws.Range("D8").Value = "'0.3561"
$ws.Range("E8").Value = "  -3.04%  "

$ws.Range("D9").Value = "'50.16"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("D10").Value = "'1.005"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").Value = "'1.209"
$ws.Range("E11").Value = "  -5.30%  "

$ws.Range("D12").Value = "'0.07950"
$ws.Range("E12").Value = "  -1.88%  "

$ws.Range("D13").Value = "'21.77"
$ws.Range("E13").Value = "  -5.95%  "

$ws.Range("D14").Value = "'6.428"
$ws.Range("E14").Value = "  -2.93%  "

$ws.Range("D15").Value = "'7.259"
$ws.Range("E15").Value = "  -4.87%  "

$ws.Range("D16").Value = "'0.00001214"
$ws.Range("E16").Value = "  -4.53%  "

$ws.Range("D17").Value = "1.579.46"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "'91.54"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").Value = "'0.06742"
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = "  -4.07%  "

$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  -3.73%  "

$ws.Range("D23").Value = "22.966.00"
$ws.Range("E23").Value = "  -1.02%  "

$ws.Range("D24").Value = "'12.66"
$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("D25").Value = "'2.361"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'2.790"
$ws.Range("E26").Value = "  -4.52%  "

$ws.Range("D27").Value = "'20.58"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").Value = "'146.94"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").Value = "'5.166"
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("D30").Value = "'131.38"
$ws.Range("E30").Value = "  -0.92%  "

$ws.Range("D31").Value = "'2.314"
$ws.Range("E31").Value = "  -4.25%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.456"
$ws.Range("E32").Value = "  -7.98%  "

$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.761.62"
$ws.Range("E33").Value = "  -1.44%  "

$ws.Range("D34").Value = "'0.9261"
$ws.Range("E34").Value = "  -6.17%  "

$ws.Range("D35").Value = "'0.07298"
$ws.Range("E35").Value = "  -5.85%  "

$ws.Range("D36").Value = "'0.02648"

$ws.Range("D37").Value = "'0.2479"
$ws.Range("E37").Value = "  -2.90%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").Value = "'9.824"
$ws.Range("E39").Value = "  -3.02%  "

$ws.Range("D40").Value = "'5.928"
$ws.Range("E40").Value = "  -6.16%  "

$ws.Range("D41").Value = "'1.324"
$ws.Range("E41").Value = "  -5.16%  "

$ws.Range("D42").Value = "'0.6791"
$ws.Range("E42").Value = "  -5.49%  "

$ws.Range("D43").Value = "'11.74"
$ws.Range("E43").Value = "  -8.34%  "

$ws.Range("D44").Value = "'14.61"
$ws.Range("E44").Value = "  -8.18%  "

$ws.Range("D45").Value = "'0.6311"
$ws.Range("E45").Value = "  -4.85%  "

$ws.Range("D46").Value = "'3.958"
$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("D47").Value = "'2.228"
$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").Value = "'130.29"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").Value = "'0.07842"
$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("D50").Value = "'1.175"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("D51").Value = "'1.162"
$ws.Range("E51").Value = "  -2.22%  "
